$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix D46: stored phone number was text, make it a real number (as in every other row) ---
$ws.Range("D46").Value = 33824565456

# --- Append new order as row 47 ---
# Start by duplicating row 46 (same client/address) so the untouched columns
# (and the blank Cantidad/Observaciones cells) keep an identical representation.
$ws.Range("A46:N46").Copy($ws.Range("A47:N47"))

# Now overwrite only the cells that differ for this new order.
$ws.Range("B47").Value = "Virginia"

# The "Fecha de Entrega" column stores plain text dates (e.g. "2025-02-27"),
# not real Excel dates. Assigning the text straight to .Value makes Excel's
# auto-detection turn it into a date serial, so instead compute it as a text
# formula result and then flatten the formula down to a plain value.
$ws.Range("E47").Formula = '="2025-02-27"'
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)

$ws.Range("F47").Value = "16:00 a 19:00"
$ws.Range("H47").Value = 20520
$ws.Range("J47").Value = "Empanadas Congeladas Pollo (12u) (x1)"
$ws.Range("N47").Value = 46

$excel.CutCopyMode = 0
